$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 322 (pushes existing rows 322-351 down to 323-352,
# matching the new dimension A1:R352).
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(322, 1).Value  = 8
$ws.Cells.Item(322, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(322, 3).Value  = "Coquimbo"
$ws.Cells.Item(322, 4).Value  = 44931
$ws.Cells.Item(322, 5).Value  = 4
$ws.Cells.Item(322, 6).Value  = 100112021
$ws.Cells.Item(322, 7).Value  = "Ají"
$ws.Cells.Item(322, 8).Value  = "Inferno"
$ws.Cells.Item(322, 9).Value  = "Primera"
$ws.Cells.Item(322, 10).Value = 480
$ws.Cells.Item(322, 11).Value = 14000
$ws.Cells.Item(322, 12).Value = 15000
$ws.Cells.Item(322, 13).Value = 14500
$ws.Cells.Item(322, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(322, 16).Value = 967
$ws.Cells.Item(322, 17).Value = 15
$ws.Cells.Item(322, 18).Value = "Hortaliza"
